$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price cells that look like plain numbers to remain text so
# they keep their exact original formatting (e.g. "82.10", "0.9977")
# instead of being auto-coerced into floating point numbers by Excel.
# (Values such as "29.270.43" have more than one "." and can never be
# parsed as a number, so they do not need this and are left alone.)
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("D13:D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20:D32").NumberFormat = "@"
$ws.Range("D34:D38").NumberFormat = "@"
$ws.Range("D40:D43").NumberFormat = "@"
$ws.Range("D45:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.270.43'
$ws.Range("E2").Value = '  -0.54%  '

# Row 3
$ws.Range("D3").Value = '1.840.82'
$ws.Range("E3").Value = '  -0.54%  '

# Row 4
$ws.Range("D4").Value = '0.9977'
$ws.Range("E4").Value = '  -0.35%  '

# Row 5
$ws.Range("D5").Value = '240.01'
$ws.Range("E5").Value = '  -0.13%  '

# Row 6
$ws.Range("D6").Value = '0.6255'
$ws.Range("E6").Value = '  -0.27%  '

# Row 7
$ws.Range("D7").Value = '0.9987'
$ws.Range("E7").Value = '  -0.27%  '

# Row 8
$ws.Range("D8").Value = '0.07475'
$ws.Range("E8").Value = '  -2.48%  '

# Row 9
$ws.Range("D9").Value = '0.2896'
$ws.Range("E9").Value = '  -0.77%  '

# Row 10
$ws.Range("D10").Value = '24.28'
$ws.Range("E10").Value = '  -2.25%  '

# Row 11
$ws.Range("D11").Value = '0.07712'
$ws.Range("E11").Value = '  -0.57%  '

# Row 12
$ws.Range("D12").Value = '1.841.21'
$ws.Range("E12").Value = '  -1.01%  '

# Row 13
$ws.Range("D13").Value = '4.983'
$ws.Range("E13").Value = '  -0.97%  '

# Row 14
$ws.Range("D14").Value = '0.6768'
$ws.Range("E14").Value = '  -0.61%  '

# Row 15
$ws.Range("D15").Value = '0.00001019'
$ws.Range("E15").Value = '  -3.73%  '

# Row 16
$ws.Range("D16").Value = '82.10'
$ws.Range("E16").Value = '  -1.63%  '

# Row 17
$ws.Range("D17").Value = '2.097.11'
$ws.Range("E17").Value = '  -0.59%  '

# Row 18
$ws.Range("D18").Value = '6.095'
$ws.Range("E18").Value = '  -1.69%  '

# Row 19
$ws.Range("D19").Value = '29.286.86'
$ws.Range("E19").Value = '  -0.61%  '

# Row 20
$ws.Range("D20").Value = '229.30'
$ws.Range("E20").Value = '  +0.28%  '

# Row 21
$ws.Range("D21").Value = '12.26'
$ws.Range("E21").Value = '  -0.52%  '

# Row 22
$ws.Range("D22").Value = '0.9987'
$ws.Range("E22").Value = '  -0.26%  '

# Row 23
$ws.Range("D23").Value = '7.370'
$ws.Range("E23").Value = '  -1.07%  '

# Row 24
$ws.Range("D24").Value = '0.9987'
$ws.Range("E24").Value = '  -0.31%  '

# Row 25
$ws.Range("D25").Value = '158.04'
$ws.Range("E25").Value = '  +0.34%  '

# Row 26
$ws.Range("D26").Value = '0.1376'
$ws.Range("E26").Value = '  -0.35%  '

# Row 27
$ws.Range("D27").Value = '8.363'
$ws.Range("E27").Value = '  -0.47%  '

# Row 28
$ws.Range("D28").Value = '17.53'
$ws.Range("E28").Value = '  -1.09%  '

# Row 29
$ws.Range("D29").Value = '1.388'
$ws.Range("E29").Value = '  +0.85%  '

# Row 30
$ws.Range("D30").Value = '1.473'
$ws.Range("E30").Value = '  +0.79%  '

# Row 31
$ws.Range("D31").Value = '0.05838'
$ws.Range("E31").Value = '  +3.95%  '

# Row 32
$ws.Range("D32").Value = '4.088'
$ws.Range("E32").Value = '  -0.85%  '

# Row 33
$ws.Range("E33").Value = '  -0.59%  '

# Row 34
$ws.Range("D34").Value = '1.811'
$ws.Range("E34").Value = '  -1.55%  '

# Row 35
$ws.Range("D35").Value = '1.142'
$ws.Range("E35").Value = '  -1.72%  '

# Row 36
$ws.Range("D36").Value = '0.6904'
$ws.Range("E36").Value = '  -1.87%  '

# Row 37
$ws.Range("D37").Value = '2.585'
$ws.Range("E37").Value = '  -0.54%  '

# Row 38
$ws.Range("D38").Value = '2.808'
$ws.Range("E38").Value = '  +2.23%  '

# Row 39
$ws.Range("D39").Value = '1.242.30'
$ws.Range("E39").Value = '  +1.73%  '

# Row 40
$ws.Range("D40").Value = '0.01816'
$ws.Range("E40").Value = '  +1.24%  '

# Row 41
$ws.Range("D41").Value = '6.519'
$ws.Range("E41").Value = '  +1.42%  '

# Row 42
$ws.Range("D42").Value = '0.9040'
$ws.Range("E42").Value = '  +0.27%  '

# Row 43
$ws.Range("D43").Value = '0.9975'
$ws.Range("E43").Value = '  -0.39%  '

# Row 44
$ws.Range("D44").Value = '1.999.23'
$ws.Range("E44").Value = '  -0.81%  '

# Row 45
$ws.Range("D45").Value = '101.33'
$ws.Range("E45").Value = '  -0.43%  '

# Row 46
$ws.Range("D46").Value = '65.79'
$ws.Range("E46").Value = '  -0.30%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.045'
$ws.Range("E47").Value = '  -1.80%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1165'
$ws.Range("E48").Value = '  +0.91%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.994'
$ws.Range("E49").Value = '  +0.13%  '

# Row 50
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").Value = '0.3930'
$ws.Range("E50").Value = '  -2.21%  '

# Row 51
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.00000000114'
$ws.Range("E51").Value = '  -5.33%  '

